$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old standalone "C3, / 100 nF" row (row 9); the row that
# follows it ("C4, C5, C6, / 100nF") shifts up to take its place.
$ws.Rows.Item(9).Delete()

# Merge C3 into the ref label of the row that now sits at row 9.
$ws.Range("A9").Value = "C3, C4, C5, C6, "

# Fix the IC1 part number / footprint labels (row 11 after the shift).
$ws.Range("C11").Value = "AP7375-50W5-7"
$ws.Range("D11").Value = "SOT95P285X140-5N"

# Leave the UI selection on the row that was edited.
$ws.Rows.Item(9).Select() | Out-Null
